# This script re-implements the commit:
# "44 advance view engine hbs concept like partials and helper"
#
# It takes the final paragraph of the document ("Axios support promise "),
# which carries the hidden "_GoBack" bookmark, and appends eight new
# paragraphs describing __dirname, Express view engines and the hbs
# (Handlebars) helper/partial APIs. The "_GoBack" bookmark is relocated
# from the old last paragraph onto the new final content paragraph, and a
# trailing empty paragraph (with the same custom tab stop) is left at the
# very end of the document, matching the author's original edit.

$d = $word.ActiveDocument

# Locate the paragraph that currently ends the document/story — this is the
# "Axios support promise " paragraph that carries the "_GoBack" bookmark.
$last = $d.Paragraphs.Last
$target = $d.Range($last.Range.Start, $last.Range.End)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00634B56" w:rsidRDefault="00634B56"><w:r><w:t xml:space="preserve">Axios support promise </w:t></w:r></w:p><w:p><w:r><w:lastRenderedPageBreak/><w:t>__dirname :- directory of your current file</w:t></w:r></w:p><w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="2535"/></w:tabs></w:pPr><w:r><w:t xml:space="preserve">View engine in express </w:t></w:r><w:r><w:t>:-</w:t></w:r><w:r><w:tab/><w:t>ejs, jade,</w:t></w:r><w:r><w:t xml:space="preserve"> handlebar</w:t></w:r><w:r><w:t>, angularjs</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>pug</w:t></w:r><w:r><w:t xml:space="preserve"> etc</w:t></w:r></w:p><w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="2535"/></w:tabs></w:pPr><w:r><w:t>Nodemon filename –e fileextension,, :- -e extenstion to which we want to watch files</w:t></w:r></w:p><w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="2535"/></w:tabs></w:pPr><w:r><w:t xml:space="preserve">HBS- Handlebars </w:t></w:r></w:p><w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="2535"/></w:tabs></w:pPr><w:r><w:t>hbs.registerPartials(__dirname+''/views/partials'')</w:t></w:r></w:p><w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="2535"/></w:tabs></w:pPr><w:r><w:t>app.set(''view engine'',''hbs'');</w:t></w:r></w:p><w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="2535"/></w:tabs></w:pPr><w:r><w:t>hbs.registerHelper(helpername,function);</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="2535"/></w:tabs></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($xml)
